$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.780.16'
$ws.Range("E2").Value = '  -7.81%  '
$ws.Range("D3").Value = '2.519.75'
$ws.Range("E3").Value = '  -3.76%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.02%  '
$ws.Range("E7").Value = '  -5.35%  '
$ws.Range("E9").Value = '  -6.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.52%  '
$ws.Range("D14").Value = '2.906.14'
$ws.Range("E14").Value = '  -3.60%  '
$ws.Range("D15").Value = '2.521.99'
$ws.Range("E15").Value = '  -3.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.863'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.00%  '
$ws.Range("E17").Value = '  -6.43%  '
$ws.Range("D18").Value = '42.837.04'
$ws.Range("E18").Value = '  -8.04%  '
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -5.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.56%  '
$ws.Range("E31").Value = '  -7.76%  '
$ws.Range("E32").Value = '  -4.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("E36").Value = '  -4.95%  '
$ws.Range("E37").Value = '  -6.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.119'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.64%  '
$ws.Range("E42").Value = '  -7.23%  '
$ws.Range("E43").Value = '  -6.18%  '
$ws.Range("D44").Value = '2.019.90'
$ws.Range("E44").Value = '  -5.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '85.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.20%  '
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("E48").Value = '  -6.74%  '
$ws.Range("D49").Value = '2.767.25'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.40%  '
$ws.Range("E51").Value = '  -7.82%  '
